$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- New test data rows (110-163): VH1-3 (high voltage), IH1-3 (high current),
# ---- and their corresponding phase-angle groups ----
$ws.Cells.Item(110,1).Value = "VH1"
$ws.Cells.Item(110,2).Value = 50
$ws.Cells.Item(110,3).Value = 49.85
$ws.Cells.Item(110,4).Value = 50.15
$ws.Cells.Item(110,5).Value = 50
$ws.Cells.Item(110,6).Formula = '=IF(AND(E110>=C110,E110<=D110), "PASS", "FAIL")'
$ws.Cells.Item(111,2).Value = 100
$ws.Cells.Item(111,3).Value = 99.7
$ws.Cells.Item(111,4).Value = 100.3
$ws.Cells.Item(111,5).Value = 100
$ws.Cells.Item(111,6).Formula = '=IF(AND(E111>=C111,E111<=D111), "PASS", "FAIL")'
$ws.Cells.Item(112,2).Value = 150
$ws.Cells.Item(112,3).Value = 149.55
$ws.Cells.Item(112,4).Value = 150.45
$ws.Cells.Item(112,5).Value = 150
$ws.Cells.Item(112,6).Formula = '=IF(AND(E112>=C112,E112<=D112), "PASS", "FAIL")'
$ws.Cells.Item(113,2).Value = 200
$ws.Cells.Item(113,3).Value = 199.4
$ws.Cells.Item(113,4).Value = 200.6
$ws.Cells.Item(113,5).Value = 200
$ws.Cells.Item(113,6).Formula = '=IF(AND(E113>=C113,E113<=D113), "PASS", "FAIL")'
$ws.Cells.Item(114,2).Value = 250
$ws.Cells.Item(114,3).Value = 249.25
$ws.Cells.Item(114,4).Value = 250.75
$ws.Cells.Item(114,5).Value = 250
$ws.Cells.Item(114,6).Formula = '=IF(AND(E114>=C114,E114<=D114), "PASS", "FAIL")'
$ws.Cells.Item(115,2).Value = 300
$ws.Cells.Item(115,3).Value = 299
$ws.Cells.Item(115,4).Value = 301
$ws.Cells.Item(115,5).Value = 301
$ws.Cells.Item(115,6).Formula = '=IF(AND(E115>=C115,E115<=D115), "PASS", "FAIL")'
$ws.Cells.Item(116,1).Value = "VH2"
$ws.Cells.Item(116,2).Value = 50
$ws.Cells.Item(116,3).Value = 49.85
$ws.Cells.Item(116,4).Value = 50.15
$ws.Cells.Item(116,5).Value = 50
$ws.Cells.Item(116,6).Formula = '=IF(AND(E116>=C116,E116<=D116), "PASS", "FAIL")'
$ws.Cells.Item(117,2).Value = 100
$ws.Cells.Item(117,3).Value = 99.7
$ws.Cells.Item(117,4).Value = 100.3
$ws.Cells.Item(117,5).Value = 100
$ws.Cells.Item(117,6).Formula = '=IF(AND(E117>=C117,E117<=D117), "PASS", "FAIL")'
$ws.Cells.Item(118,2).Value = 150
$ws.Cells.Item(118,3).Value = 149.55
$ws.Cells.Item(118,4).Value = 150.45
$ws.Cells.Item(118,5).Value = 150
$ws.Cells.Item(118,6).Formula = '=IF(AND(E118>=C118,E118<=D118), "PASS", "FAIL")'
$ws.Cells.Item(119,2).Value = 200
$ws.Cells.Item(119,3).Value = 199.4
$ws.Cells.Item(119,4).Value = 200.6
$ws.Cells.Item(119,5).Value = 200
$ws.Cells.Item(119,6).Formula = '=IF(AND(E119>=C119,E119<=D119), "PASS", "FAIL")'
$ws.Cells.Item(120,2).Value = 250
$ws.Cells.Item(120,3).Value = 249.25
$ws.Cells.Item(120,4).Value = 250.75
$ws.Cells.Item(120,5).Value = 250
$ws.Cells.Item(120,6).Formula = '=IF(AND(E120>=C120,E120<=D120), "PASS", "FAIL")'
$ws.Cells.Item(121,2).Value = 300
$ws.Cells.Item(121,3).Value = 299
$ws.Cells.Item(121,4).Value = 301
$ws.Cells.Item(121,5).Value = 301
$ws.Cells.Item(121,6).Formula = '=IF(AND(E121>=C121,E121<=D121), "PASS", "FAIL")'
$ws.Cells.Item(122,1).Value = "VH3"
$ws.Cells.Item(122,2).Value = 50
$ws.Cells.Item(122,3).Value = 49.85
$ws.Cells.Item(122,4).Value = 50.15
$ws.Cells.Item(122,5).Value = 50
$ws.Cells.Item(122,6).Formula = '=IF(AND(E122>=C122,E122<=D122), "PASS", "FAIL")'
$ws.Cells.Item(123,2).Value = 100
$ws.Cells.Item(123,3).Value = 99.7
$ws.Cells.Item(123,4).Value = 100.3
$ws.Cells.Item(123,5).Value = 100
$ws.Cells.Item(123,6).Formula = '=IF(AND(E123>=C123,E123<=D123), "PASS", "FAIL")'
$ws.Cells.Item(124,2).Value = 150
$ws.Cells.Item(124,3).Value = 149.55
$ws.Cells.Item(124,4).Value = 150.45
$ws.Cells.Item(124,5).Value = 150
$ws.Cells.Item(124,6).Formula = '=IF(AND(E124>=C124,E124<=D124), "PASS", "FAIL")'
$ws.Cells.Item(125,2).Value = 200
$ws.Cells.Item(125,3).Value = 199.4
$ws.Cells.Item(125,4).Value = 200.6
$ws.Cells.Item(125,5).Value = 200
$ws.Cells.Item(125,6).Formula = '=IF(AND(E125>=C125,E125<=D125), "PASS", "FAIL")'
$ws.Cells.Item(126,2).Value = 250
$ws.Cells.Item(126,3).Value = 249.25
$ws.Cells.Item(126,4).Value = 250.75
$ws.Cells.Item(126,5).Value = 250
$ws.Cells.Item(126,6).Formula = '=IF(AND(E126>=C126,E126<=D126), "PASS", "FAIL")'
$ws.Cells.Item(127,2).Value = 300
$ws.Cells.Item(127,3).Value = 299
$ws.Cells.Item(127,4).Value = 301
$ws.Cells.Item(127,5).Value = 301
$ws.Cells.Item(127,6).Formula = '=IF(AND(E127>=C127,E127<=D127), "PASS", "FAIL")'
$ws.Cells.Item(128,1).Value = "IH1"
$ws.Cells.Item(128,2).Value = 1
$ws.Cells.Item(128,3).Value = 0.997
$ws.Cells.Item(128,4).Value = 1.003
$ws.Cells.Item(128,5).Value = 0.99952531
$ws.Cells.Item(128,6).Formula = '=IF(AND(E128>=C128,E128<=D128), "PASS", "FAIL")'
$ws.Cells.Item(129,2).Value = 2
$ws.Cells.Item(129,3).Value = 1.994
$ws.Cells.Item(129,4).Value = 2.006
$ws.Cells.Item(129,5).Value = 1.99924719
$ws.Cells.Item(129,6).Formula = '=IF(AND(E129>=C129,E129<=D129), "PASS", "FAIL")'
$ws.Cells.Item(130,2).Value = 3
$ws.Cells.Item(130,3).Value = 2.991
$ws.Cells.Item(130,4).Value = 3.009
$ws.Cells.Item(130,5).Value = 2.9996326
$ws.Cells.Item(130,6).Formula = '=IF(AND(E130>=C130,E130<=D130), "PASS", "FAIL")'
$ws.Cells.Item(131,2).Value = 4
$ws.Cells.Item(131,3).Value = 3.988
$ws.Cells.Item(131,4).Value = 4.012
$ws.Cells.Item(131,5).Value = 3.99974108
$ws.Cells.Item(131,6).Formula = '=IF(AND(E131>=C131,E131<=D131), "PASS", "FAIL")'
$ws.Cells.Item(132,2).Value = 5
$ws.Cells.Item(132,3).Value = 4.985
$ws.Cells.Item(132,4).Value = 5.015
$ws.Cells.Item(132,5).Value = 5.00013781
$ws.Cells.Item(132,6).Formula = '=IF(AND(E132>=C132,E132<=D132), "PASS", "FAIL")'
$ws.Cells.Item(133,2).Value = 6
$ws.Cells.Item(133,3).Value = 5.982
$ws.Cells.Item(133,4).Value = 6.018
$ws.Cells.Item(133,5).Value = 6.00014973
$ws.Cells.Item(133,6).Formula = '=IF(AND(E133>=C133,E133<=D133), "PASS", "FAIL")'
$ws.Cells.Item(134,1).Value = "IH2"
$ws.Cells.Item(134,2).Value = 1
$ws.Cells.Item(134,3).Value = 0.997
$ws.Cells.Item(134,4).Value = 1.003
$ws.Cells.Item(134,5).Value = 0.9997226
$ws.Cells.Item(134,6).Formula = '=IF(AND(E134>=C134,E134<=D134), "PASS", "FAIL")'
$ws.Cells.Item(135,2).Value = 2
$ws.Cells.Item(135,3).Value = 1.994
$ws.Cells.Item(135,4).Value = 2.006
$ws.Cells.Item(135,5).Value = 1.9996196
$ws.Cells.Item(135,6).Formula = '=IF(AND(E135>=C135,E135<=D135), "PASS", "FAIL")'
$ws.Cells.Item(136,2).Value = 3
$ws.Cells.Item(136,3).Value = 2.991
$ws.Cells.Item(136,4).Value = 3.009
$ws.Cells.Item(136,5).Value = 3.00004983
$ws.Cells.Item(136,6).Formula = '=IF(AND(E136>=C136,E136<=D136), "PASS", "FAIL")'
$ws.Cells.Item(137,2).Value = 4
$ws.Cells.Item(137,3).Value = 3.988
$ws.Cells.Item(137,4).Value = 4.012
$ws.Cells.Item(137,5).Value = 4.0004735
$ws.Cells.Item(137,6).Formula = '=IF(AND(E137>=C137,E137<=D137), "PASS", "FAIL")'
$ws.Cells.Item(138,2).Value = 5
$ws.Cells.Item(138,3).Value = 4.985
$ws.Cells.Item(138,4).Value = 5.015
$ws.Cells.Item(138,5).Value = 5.00082779
$ws.Cells.Item(138,6).Formula = '=IF(AND(E138>=C138,E138<=D138), "PASS", "FAIL")'
$ws.Cells.Item(139,2).Value = 6
$ws.Cells.Item(139,3).Value = 5.982
$ws.Cells.Item(139,4).Value = 6.018
$ws.Cells.Item(139,5).Value = 6.00077391
$ws.Cells.Item(139,6).Formula = '=IF(AND(E139>=C139,E139<=D139), "PASS", "FAIL")'
$ws.Cells.Item(140,1).Value = "IH3"
$ws.Cells.Item(140,2).Value = 1
$ws.Cells.Item(140,3).Value = 0.997
$ws.Cells.Item(140,4).Value = 1.003
$ws.Cells.Item(140,5).Value = 0.99990565
$ws.Cells.Item(140,6).Formula = '=IF(AND(E140>=C140,E140<=D140), "PASS", "FAIL")'
$ws.Cells.Item(141,2).Value = 2
$ws.Cells.Item(141,3).Value = 1.994
$ws.Cells.Item(141,4).Value = 2.006
$ws.Cells.Item(141,5).Value = 1.99978471
$ws.Cells.Item(141,6).Formula = '=IF(AND(E141>=C141,E141<=D141), "PASS", "FAIL")'
$ws.Cells.Item(142,2).Value = 3
$ws.Cells.Item(142,3).Value = 2.991
$ws.Cells.Item(142,4).Value = 3.009
$ws.Cells.Item(142,5).Value = 3.00055504
$ws.Cells.Item(142,6).Formula = '=IF(AND(E142>=C142,E142<=D142), "PASS", "FAIL")'
$ws.Cells.Item(143,2).Value = 4
$ws.Cells.Item(143,3).Value = 3.988
$ws.Cells.Item(143,4).Value = 4.012
$ws.Cells.Item(143,5).Value = 4.00096369
$ws.Cells.Item(143,6).Formula = '=IF(AND(E143>=C143,E143<=D143), "PASS", "FAIL")'
$ws.Cells.Item(144,2).Value = 5
$ws.Cells.Item(144,3).Value = 4.985
$ws.Cells.Item(144,4).Value = 5.015
$ws.Cells.Item(144,5).Value = 5.00128889
$ws.Cells.Item(144,6).Formula = '=IF(AND(E144>=C144,E144<=D144), "PASS", "FAIL")'
$ws.Cells.Item(145,2).Value = 6
$ws.Cells.Item(145,3).Value = 5.982
$ws.Cells.Item(145,4).Value = 6.018
$ws.Cells.Item(145,5).Value = 6.00170803
$ws.Cells.Item(145,6).Formula = '=IF(AND(E145>=C145,E145<=D145), "PASS", "FAIL")'
$ws.Cells.Item(146,1).Value = "Phase(VH1)"
$ws.Cells.Item(146,2).Value = 60
$ws.Cells.Item(146,3).Value = 59.75
$ws.Cells.Item(146,4).Value = 60.25
$ws.Cells.Item(146,5).Value = 60.26304637721221
$ws.Cells.Item(146,6).Formula = '=IF(AND(E146>=C146,E146<=D146), "PASS", "FAIL")'
$ws.Cells.Item(147,2).Value = 120
$ws.Cells.Item(147,3).Value = 119.75
$ws.Cells.Item(147,4).Value = 120.25
$ws.Cells.Item(147,5).Value = 120.2972119614733
$ws.Cells.Item(147,6).Formula = '=IF(AND(E147>=C147,E147<=D147), "PASS", "FAIL")'
$ws.Cells.Item(148,2).Value = 180
$ws.Cells.Item(148,3).Value = 179.75
$ws.Cells.Item(148,4).Value = 180.25
$ws.Cells.Item(148,5).Value = 179.7593575491723
$ws.Cells.Item(148,6).Formula = '=IF(AND(E148>=C148,E148<=D148), "PASS", "FAIL")'
$ws.Cells.Item(149,1).Value = "Phase(VH2)"
$ws.Cells.Item(149,2).Value = 60
$ws.Cells.Item(149,3).Value = 59.75
$ws.Cells.Item(149,4).Value = 60.25
$ws.Cells.Item(149,5).Value = 60.23295852263566
$ws.Cells.Item(149,6).Formula = '=IF(AND(E149>=C149,E149<=D149), "PASS", "FAIL")'
$ws.Cells.Item(150,2).Value = 120
$ws.Cells.Item(150,3).Value = 119.75
$ws.Cells.Item(150,4).Value = 120.25
$ws.Cells.Item(150,5).Value = 120.2331134184667
$ws.Cells.Item(150,6).Formula = '=IF(AND(E150>=C150,E150<=D150), "PASS", "FAIL")'
$ws.Cells.Item(151,2).Value = 180
$ws.Cells.Item(151,3).Value = 179.75
$ws.Cells.Item(151,4).Value = 180.25
$ws.Cells.Item(151,5).Value = 179.768394454297
$ws.Cells.Item(151,6).Formula = '=IF(AND(E151>=C151,E151<=D151), "PASS", "FAIL")'
$ws.Cells.Item(152,1).Value = "Phase(VH3)"
$ws.Cells.Item(152,2).Value = 60
$ws.Cells.Item(152,3).Value = 59.75
$ws.Cells.Item(152,4).Value = 60.25
$ws.Cells.Item(152,5).Value = 60.21943581112865
$ws.Cells.Item(152,6).Formula = '=IF(AND(E152>=C152,E152<=D152), "PASS", "FAIL")'
$ws.Cells.Item(153,2).Value = 120
$ws.Cells.Item(153,3).Value = 119.75
$ws.Cells.Item(153,4).Value = 120.25
$ws.Cells.Item(153,5).Value = 120.2084549933066
$ws.Cells.Item(153,6).Formula = '=IF(AND(E153>=C153,E153<=D153), "PASS", "FAIL")'
$ws.Cells.Item(154,2).Value = 180
$ws.Cells.Item(154,3).Value = 179.75
$ws.Cells.Item(154,4).Value = 180.25
$ws.Cells.Item(154,5).Value = 179.8082515568856
$ws.Cells.Item(154,6).Formula = '=IF(AND(E154>=C154,E154<=D154), "PASS", "FAIL")'
$ws.Cells.Item(155,1).Value = "Phase(IH1)"
$ws.Cells.Item(155,2).Value = 60
$ws.Cells.Item(155,3).Value = 59.75
$ws.Cells.Item(155,4).Value = 60.25
$ws.Cells.Item(155,5).Value = 60.2292403880278
$ws.Cells.Item(155,6).Formula = '=IF(AND(E155>=C155,E155<=D155), "PASS", "FAIL")'
$ws.Cells.Item(156,2).Value = 120
$ws.Cells.Item(156,3).Value = 119.75
$ws.Cells.Item(156,4).Value = 120.25
$ws.Cells.Item(156,5).Value = 120.2336267040371
$ws.Cells.Item(156,6).Formula = '=IF(AND(E156>=C156,E156<=D156), "PASS", "FAIL")'
$ws.Cells.Item(157,2).Value = 180
$ws.Cells.Item(157,3).Value = 179.75
$ws.Cells.Item(157,4).Value = 180.25
$ws.Cells.Item(157,5).Value = 179.7770611051297
$ws.Cells.Item(157,6).Formula = '=IF(AND(E157>=C157,E157<=D157), "PASS", "FAIL")'
$ws.Cells.Item(158,1).Value = "Phase(IH2)"
$ws.Cells.Item(158,2).Value = 60
$ws.Cells.Item(158,3).Value = 59.75
$ws.Cells.Item(158,4).Value = 60.25
$ws.Cells.Item(158,5).Value = 59.80549734096019
$ws.Cells.Item(158,6).Formula = '=IF(AND(E158>=C158,E158<=D158), "PASS", "FAIL")'
$ws.Cells.Item(159,2).Value = 120
$ws.Cells.Item(159,3).Value = 119.75
$ws.Cells.Item(159,4).Value = 120.25
$ws.Cells.Item(159,5).Value = 119.7908741234901
$ws.Cells.Item(159,6).Formula = '=IF(AND(E159>=C159,E159<=D159), "PASS", "FAIL")'
$ws.Cells.Item(160,2).Value = 180
$ws.Cells.Item(160,3).Value = 179.75
$ws.Cells.Item(160,4).Value = 180.25
$ws.Cells.Item(160,5).Value = 179.8241476184619
$ws.Cells.Item(160,6).Formula = '=IF(AND(E160>=C160,E160<=D160), "PASS", "FAIL")'
$ws.Cells.Item(161,1).Value = "Phase(IH3)"
$ws.Cells.Item(161,2).Value = 60
$ws.Cells.Item(161,3).Value = 59.75
$ws.Cells.Item(161,4).Value = 60.25
$ws.Cells.Item(161,5).Value = 59.78603605450781
$ws.Cells.Item(161,6).Formula = '=IF(AND(E161>=C161,E161<=D161), "PASS", "FAIL")'
$ws.Cells.Item(162,2).Value = 120
$ws.Cells.Item(162,3).Value = 119.75
$ws.Cells.Item(162,4).Value = 120.25
$ws.Cells.Item(162,5).Value = 119.7979768913345
$ws.Cells.Item(162,6).Formula = '=IF(AND(E162>=C162,E162<=D162), "PASS", "FAIL")'
$ws.Cells.Item(163,2).Value = 180
$ws.Cells.Item(163,3).Value = 179.75
$ws.Cells.Item(163,4).Value = 180.25
$ws.Cells.Item(163,5).Value = 179.8031824372122
$ws.Cells.Item(163,6).Formula = '=IF(AND(E163>=C163,E163<=D163), "PASS", "FAIL")'

# ---- Conditional formatting ("PASS"/"FAIL" highlighting) ----
# Existing rule's range grows to cover the first new block (F2:F118) -
# extend in place so the rule keeps its original dxf (colours unchanged).
$existingRange = $ws.Range("F2:F109")
$existingFcs = $existingRange.FormatConditions
for ($i = 1; $i -le $existingFcs.Count; $i++) {
    $existingFcs.Item($i).ModifyAppliesToRange($ws.Range("F2:F118")) | Out-Null
}

$greenColor = 5287936   # RGB(0,176,80)  -> FF00B050
$redColor   = 255       # RGB(255,0,0)   -> FFFF0000

function Add-PassFailCF($rangeAddress) {
    $rng = $ws.Range($rangeAddress)
    $failRule = $rng.FormatConditions.Add(9, 0, $null, $null, "FAIL")
    $failRule.Interior.Color = $redColor
    $passRule = $rng.FormatConditions.Add(9, 0, $null, $null, "PASS")
    $passRule.Interior.Color = $greenColor
}

Add-PassFailCF "F119:F127"
Add-PassFailCF "F128:F145"
Add-PassFailCF "F146:F154"
Add-PassFailCF "F155:F163"

# ---- View state: scroll position & selection left by the editor ----
$ws.Range("B122:D127").Select()

Write-Host "edit complete"
